$wb = $excel.ActiveWorkbook

$wsLogin  = $wb.Worksheets.Item("Login")
$wsLogOff = $wb.Worksheets.Item("LogOff")

# --- LogOff sheet ("tech@siigo.copm" test row, mistyped address used as a
#     negative / listener test case) ---
$wsLogOff.Range("A3").Value = "tech@siigo.copm"
$wsLogOff.Range("B3").Value = 1111
$wsLogOff.Hyperlinks.Add($wsLogOff.Range("A3"), "mailto:tech@siigo.copm")
$wsLogOff.Range("A3").Style = "Hipervínculo"

# --- Login sheet (corrected "tech@siigo.com" test row) ---
$wsLogin.Range("B3").Value = "tech@siigo.com"
$wsLogin.Range("C3").Value = 2222
$wsLogin.Hyperlinks.Add($wsLogin.Range("B3"), "mailto:tech@siigo.com")
$wsLogin.Range("B3").Style = "Hipervínculo"

# --- Selection / active-tab bookkeeping left behind by running the test
#     case bindings: LogOff's cursor resets to C13, and Login (now bound to
#     the test listeners) becomes the active sheet, selected at C9. ---
$wsLogOff.Range("C13").Select() | Out-Null
$wsLogin.Activate()
$wsLogin.Range("C9").Select() | Out-Null
